# Rewrite the player roster table: reorder rows, update team/position
# assignments, and drop the "Bennedict Mathurin" / "Indiana Pacers" row so
# the table now has 16 players (rows 2-17) instead of 17 (rows 2-18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Tyler Herro",     "PG,SG",    "Miami Heat"),
    @("Josh Giddey",     "PG,SG,SF", "Chicago Bulls"),
    @("Caris LeVert",    "SG,SF",    "Cleveland Cavaliers"),
    @("Miles Bridges",   "SF,PF",    "Charlotte Hornets"),
    @("DeMar DeRozan",   "SF,PF",    "Sacramento Kings"),
    @("Nikola Vucevic",  "PF,C",     "Chicago Bulls"),
    @("Evan Mobley",     "PF,C",     "Cleveland Cavaliers"),
    @("Nick Richards",   "C",        "Charlotte Hornets"),
    @("De'Aaron Fox",    "PG",       "Sacramento Kings"),
    @("Gradey Dick",     "SG,SF",    "Toronto Raptors"),
    @("Brook Lopez",     "C",        "Milwaukee Bucks"),
    @("Santi Aldama",    "PF,C",     "Memphis Grizzlies"),
    @("Luka Doncic",     "PG,SG",    "Dallas Mavericks"),
    @("Ja Morant",       "PG",       "Memphis Grizzlies"),
    @("Mikal Bridges",   "SG,SF,PF", "New York Knicks"),
    @("Scottie Barnes",  "SG,SF,PF", "Toronto Raptors")
)

# Clear out the old table body (rows 2 through 18) before writing the new,
# shorter table (rows 2 through 17).
$ws.Range("A2:C18").ClearContents()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
